$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    0.7724775713279004,
    0.732260308402072,
    0.2534333533236305,
    0.1276061961020767,
    0.8962634008536721,
    0.7739018424399541,
    0.1223483126088004,
    0.2652424226402139,
    1.000150259811752,
    1.728893897005647,
    0.9377931549521943,
    1.633505894680542,
    1.7764735337147
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $values[$i]
}
